# Auto-generated edit script: update column F ('想去人数' / want-to-go count)
# values across all 4 sheets, plus one column G ('最低票价') correction on
# sheet 4 row 10, per the upstream data refresh commit.

$wb = $excel.ActiveWorkbook

# Sheet 1: '展览' (by position, Chinese sheet names print oddly
# through this harness's stdout but addressing by name also works;
# positional Item() is used here to keep things unambiguous).
$ws1 = $wb.Worksheets.Item(1)
$fUpdates1 = @(
    @(3, 257),
    @(4, 573),
    @(5, 2523),
    @(7, 158),
    @(9, 241),
    @(10, 5111),
    @(11, 6248),
    @(12, 89),
    @(13, 1429),
    @(14, 1361),
    @(15, 579),
    @(16, 6880),
    @(18, 35),
    @(19, 3),
    @(20, 65),
    @(21, 4624),
    @(22, 377),
    @(23, 57),
    @(24, 2268),
    @(25, 1239),
    @(26, 426),
    @(27, 1141),
    @(28, 201),
    @(29, 89),
    @(30, 68),
    @(31, 141),
    @(32, 363),
    @(33, 1258),
    @(34, 1973),
    @(35, 213),
    @(36, 490),
    @(38, 1345),
    @(39, 586),
    @(42, 1084),
    @(43, 2371),
    @(44, 39),
    @(45, 43),
    @(46, 91),
    @(48, 66)
)
foreach ($pair in $fUpdates1) {
    $ws1.Cells.Item($pair[0], 6).Value = $pair[1]
}

# Sheet 2: '演出' (by position, Chinese sheet names print oddly
# through this harness's stdout but addressing by name also works;
# positional Item() is used here to keep things unambiguous).
$ws2 = $wb.Worksheets.Item(2)
$fUpdates2 = @(
    @(4, 426),
    @(7, 131),
    @(12, 370),
    @(13, 255),
    @(14, 160),
    @(15, 29),
    @(23, 140),
    @(24, 31),
    @(28, 265),
    @(35, 1),
    @(36, 3),
    @(38, 5)
)
foreach ($pair in $fUpdates2) {
    $ws2.Cells.Item($pair[0], 6).Value = $pair[1]
}

# Sheet 3: '本地生活' (by position, Chinese sheet names print oddly
# through this harness's stdout but addressing by name also works;
# positional Item() is used here to keep things unambiguous).
$ws3 = $wb.Worksheets.Item(3)
$fUpdates3 = @(
    @(4, 490),
    @(6, 1639),
    @(7, 532),
    @(8, 1241),
    @(10, 1715),
    @(11, 2090),
    @(12, 555),
    @(13, 460)
)
foreach ($pair in $fUpdates3) {
    $ws3.Cells.Item($pair[0], 6).Value = $pair[1]
}

# Sheet 4: '全部类型' (by position, Chinese sheet names print oddly
# through this harness's stdout but addressing by name also works;
# positional Item() is used here to keep things unambiguous).
$ws4 = $wb.Worksheets.Item(4)
$fUpdates4 = @(
    @(2, 490),
    @(3, 1639),
    @(4, 257),
    @(5, 573),
    @(6, 2523),
    @(7, 158),
    @(8, 1241),
    @(9, 241),
    @(10, 5111),
    @(11, 555),
    @(14, 89),
    @(16, 1429),
    @(17, 1361),
    @(18, 579),
    @(19, 6880),
    @(21, 460),
    @(24, 35),
    @(25, 4624),
    @(26, 2268),
    @(27, 426),
    @(28, 1141),
    @(29, 201),
    @(30, 89),
    @(31, 68),
    @(32, 255),
    @(33, 141),
    @(34, 363),
    @(35, 1258),
    @(36, 1973),
    @(37, 213),
    @(38, 490),
    @(41, 1345),
    @(42, 140),
    @(46, 1084),
    @(47, 2371),
    @(48, 66),
    @(49, 5)
)
foreach ($pair in $fUpdates4) {
    $ws4.Cells.Item($pair[0], 6).Value = $pair[1]
}

# Sheet 4 ('全部类型') row 10 also has a column-G ('最低票价') change
# alongside its column-F change above.
$ws4.Cells.Item(10, 7).Value = 69
